$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

# The import picked up a new (more recent) contract-note row. Shift the
# existing data rows (5-19) down by one to make room for it at row 5,
# then fill row 5 with the newly detected transaction.

$dateFmt = "yyyy-mm-dd h:mm:ss"

for ($r = 19; $r -ge 5; $r--) {
    $dst = $r + 1
    $ws.Range("A$dst`:J$dst").ClearContents()

    foreach ($col in @("A","B","C","D","E","F","G","H","I")) {
        $srcCell = $ws.Range("$col$r")
        $v = $srcCell.Value2()
        if ($v -ne "") {
            $ws.Range("$col$dst").Value = $v
        }
    }

    $ws.Range("A$dst").NumberFormat = $dateFmt
    $ws.Range("J$dst").Formula = "=Index!`$C`$2"
}

# New row 5: the duplicate-detection import of the latest contract note.
$ws.Range("A5:J5").ClearContents()

$ws.Range("A5").Value = 46063
$ws.Range("A5").NumberFormat = $dateFmt
$ws.Range("B5").Value = "BSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2022.1
$ws.Range("F5").Value = 4072.87
$ws.Range("G5").Value = "CN#252611730667"
$ws.Range("H5").Value = 4.03
$ws.Range("I5").Value = 24.64
$ws.Range("J5").Formula = "=Index!`$C`$2"
